$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff shows the 4 data rows (2-5) being rotated: new row2 gets old row4's
# values, new row3 gets old row5's values, new row4 gets old row3's values,
# and new row5 gets old row2's values (columns D, M, N, O, P, S only).
# Capture the original values first, then write the new arrangement.

$orig = @{}
foreach ($r in 2..5) {
    $orig[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2   # Fecha
        M = $ws.Cells.Item($r, 13).Value2  # Volumen
        N = $ws.Cells.Item($r, 14).Value2  # Precio minimo
        O = $ws.Cells.Item($r, 15).Value2  # Precio maximo
        P = $ws.Cells.Item($r, 16).Value2  # Precio promedio ponderado
        S = $ws.Cells.Item($r, 19).Value2  # Precio $/Kg
    }
}

# mapping: target row -> source row (using original pre-edit values)
$mapping = @{ 2 = 4; 3 = 5; 4 = 3; 5 = 2 }

foreach ($target in $mapping.Keys) {
    $source = $mapping[$target]
    $vals = $orig[$source]
    $ws.Cells.Item($target, 4).Value = $vals.D
    $ws.Cells.Item($target, 13).Value = $vals.M
    $ws.Cells.Item($target, 14).Value = $vals.N
    $ws.Cells.Item($target, 15).Value = $vals.O
    $ws.Cells.Item($target, 16).Value = $vals.P
    $ws.Cells.Item($target, 19).Value = $vals.S
}
